# Add "ApPredict version information" worksheet with the build/version
# metadata for the ApPredict run, at the end of the workbook, and
# re-write the (unchanged, just re-serialised) value in
# "Input Values"!B12 so the workbook matches the re-saved file.

$wb = $excel.ActiveWorkbook

# --- New worksheet -------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([Type]::Missing, $lastSheet)
$ws.Name = "ApPredict version information"

# Two-column (label, value) table describing the ApPredict / Chaste build
# used to run the simulation.
$data = @(
    @("ApPredict Version", "37cc5a6"),
    @("Chaste Version", "2019.1.682dce0"),
    @("Modified", $true),
    @("Build options", "GccOpt, shared libraries"),
    @("OS info", "Linux d09b088bdc9f 4.15.0-161-generic #169-Ubuntu SMP Fri Oct 15 13:41:54 UTC 2021 x86_64"),
    @("Compiler", "gcc, version b'9.3.0'"),
    @("Compiler flags", "-O3 -std=c++14"),
    @("XSD", "4.0.0"),
    @("VTK", "no"),
    @("Xerces", "3.2.0"),
    @("SUNDIALS", "2.5.0"),
    @("HDF5", "1.8.16"),
    @("Boost", "1.65.1"),
    @("PETSc", "3.12.4"),
    @("Parmetis", "4.0.3"),
    @("Ap Predict arguments", " --pacing-freq 1 --pacing-max-time 5 --plasma-conc-high 100 --plasma-conc-low 0 --plasma-conc-count 4 --plasma-conc-logscale true --model 1")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $labelCell = $ws.Cells.Item($row, 1)
    $valueCell = $ws.Cells.Item($row, 2)

    $labelCell.Value = $data[$i][0]

    $value = $data[$i][1]
    if ($value -is [bool]) {
        $valueCell.Value = $value
    } else {
        # Force text interpretation so version strings that look numeric
        # or date-like (e.g. "1.8.16") are stored verbatim instead of
        # being auto-converted (e.g. to a date serial). ClearFormats()
        # afterwards drops the temporary "@" number format again so the
        # cell is left with the regular default style.
        $valueCell.NumberFormat = "@"
        $valueCell.Value = $value
        $valueCell.ClearFormats()
    }
}

$ws.Range("A1:B16").Select() | Out-Null

# --- Unrelated re-save touch on the first sheet --------------------------
# The original workbook recorded 44.716 here; keep the same numeric value.
$input = $wb.Worksheets.Item("Input Values")
$input.Cells.Item(12, 2).Value = 44.716

Write-Output "Added 'ApPredict version information' sheet"
